# edit.ps1 - applies the PONG GAME assignment report edits:
#  - strike through completed/obsolete "difficulty" and "game management" list items
#    (keeping the leading "N, " numbering unstruck)
#  - clean up the "deliverables" sentence (merge fragmented runs, drop proofErr marks)
#  - strike through the whole "Report should contain..." section headings/body
#    (Cover page, Introduction, Context/Rationale, Specification, Design, Build,
#     Test, Conclusion) since that work is now complete
#  - add a "_Hlk94803562" bookmark around "Context/Rationale"
#  - remove the stray "_GoBack" bookmark paragraph

$d = $word.ActiveDocument

function Strike-Range([object]$rng) {
    $rng.Font.StrikeThrough = 1
}

function Find-And-Strike([string]$text) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $rng.Font.StrikeThrough = 1
    }
    return $ok
}

function Get-Paragraph-Containing([string]$needle) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. "3, add another ball" -> "3, " (normal) + "add another ball" (struck)
# ---------------------------------------------------------------------------
$p = Get-Paragraph-Containing("add another ball")
$p.Range.Font.StrikeThrough = 1
Find-And-Strike("add another ball") | Out-Null
$rng = $d.Content
$rng.Find.Execute("3, ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 0

# ---------------------------------------------------------------------------
# 2. "4, increase speed (reduce delay)" -> "4, " (normal) + rest (struck)
# ---------------------------------------------------------------------------
$p = Get-Paragraph-Containing("increase speed (reduce delay)")
$p.Range.Font.StrikeThrough = 1
$rng = $d.Content
$rng.Find.Execute("4, ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 0

# ---------------------------------------------------------------------------
# 3. "1, restart, wait for button press" -> "1, " (normal) + rest (struck)
# ---------------------------------------------------------------------------
$p = Get-Paragraph-Containing("restart, wait for button press")
$p.Range.Font.StrikeThrough = 1
$rng = $d.Content
$rng.Find.Execute("1, ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 0

# ---------------------------------------------------------------------------
# 4. "2, game over, lives count, high score" -> "2, " stays normal, rest struck
#    (no paragraph-mark strike here, only the text after "2, ")
# ---------------------------------------------------------------------------
Find-And-Strike("game over, lives count, high score") | Out-Null

# ---------------------------------------------------------------------------
# 5. Clean up the deliverables sentence (merge runs, drop proofErr wrappers)
# ---------------------------------------------------------------------------
$rng = $d.Content
$sentence = "The deliverables for this assignment are a working demonstration of  the game along with a written report."
$rng.Find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, $sentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Strike through the whole "report sections" block (now finished/obsolete)
# ---------------------------------------------------------------------------
foreach ($needle in @(
        "Cover page",
        "Introduction:",
        "Context/Rationale",
        "Specification:",
        "Design",
        "Build",
        "Test:",
        "Conclusion"
    )) {
    $p = Get-Paragraph-Containing($needle)
    if ($p -ne $null) {
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------------
# 7. Add the "_Hlk94803562" bookmark around "Context/Rationale"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Context/Rationale", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_Hlk94803562", $rng) | Out-Null

# ---------------------------------------------------------------------------
# 8. Remove the stray "_GoBack" bookmark (leaves a bare empty paragraph)
# ---------------------------------------------------------------------------
$gb = $d.Bookmarks.Item("_GoBack")
if ($gb -ne $null) {
    $gb.Delete()
}
